$wb = $excel.ActiveWorkbook

# --- "Range Status" sheet: zero out column B values, clear column C entirely ---
$wsRange = $wb.Worksheets.Item("Range Status")
for ($r = 2; $r -le 7; $r++) {
    $wsRange.Cells.Item($r, 2).Value = 0
    $wsRange.Cells.Item($r, 3).ClearContents()
}

# --- "Species qualification" sheet: Range Analysis row B5 220 -> 0 ---
$wsSpecies = $wb.Worksheets.Item("Species qualification")
$wsSpecies.Range("B5").Value = 0

# --- "High Priority break-up" sheet: add new-species columns D2/E2 ---
$wsBreakup = $wb.Worksheets.Item("High Priority break-up")
$wsBreakup.Range("D2").Value = 1
$wsBreakup.Range("E2").Value = 100
